# Apply updated loading_percent values (case with 380 kV) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns that hold data values (B,C,E,F,G,J,K,M,O); D,H,I,L,N stay 0, column A is the index.
$cols = @(2, 3, 5, 6, 7, 10, 11, 13, 15)

# New values for rows 2..25, one row per array element, values ordered to match $cols.
$newValues = @(
    @(8.970205258459627, 5.305302056174891, 20.27290231388671, 41.26409604810227, 3.656083476421752, 8.483964644306925, 8.253613472711759, 17.46786279019956, 24.35991104550933),
    @(8.703134417172292, 5.208580668762067, 20.15962001895475, 41.22998143551515, 3.657846231962177, 8.500914905619117, 8.074554752496875, 17.3429975434606, 24.46385021059953),
    @(8.536299678792405, 5.147851050857881, 20.09387016523802, 41.21945345361361, 3.658985186785564, 8.511844397578335, 7.963646665297169, 17.26883611297267, 24.53275304325943),
    @(8.467701973254696, 5.122786493359867, 20.06805464939224, 41.21778639649991, 3.65946360144546, 8.516429907700921, 7.918269542494857, 17.23927047784229, 24.56210786271533),
    @(8.45627770060903, 5.118606035743578, 20.06382766290533, 41.21766807505197, 3.659543905710552, 8.517199291515368, 7.910725603788659, 17.2344014567702, 24.56705923878046),
    @(8.535376875486175, 5.147514275883328, 20.09351802123033, 41.21942034705409, 3.658991580969294, 8.511905705780318, 7.963035344224521, 17.26843469247461, 24.5331437673585),
    @(8.878772537153599, 5.27224331985339, 20.23306572004115, 41.25017357288316, 3.656679550959674, 8.48970100351602, 8.192110616847103, 17.42430468281679, 24.39469293255363),
    @(9.525259204136855, 5.505318464939883, 20.53592694647902, 41.39291366257164, 3.652592837963502, 8.450280711537008, 8.63108286592559, 17.74868497682676, 24.16361803971573),
    @(9.97863659699027, 5.668443032946427, 20.77483344334064, 41.54759999950593, 3.649860060537856, 8.423806061380578, 8.944083439854834, 17.99679906682165, 24.01861498029422),
    @(10.17925722581177, 5.740678666954232, 20.886753756367, 41.62865103083018, 3.648674806391598, 8.412296616681623, 9.083779707536575, 18.11144020608674, 23.95805409640557),
    @(10.25434997066928, 5.767732714325401, 20.92957228565679, 41.66086418680336, 3.648234261158259, 8.408014661771041, 9.136244251371473, 18.15507588993456, 23.93590035477809),
    @(10.23821747934623, 5.76191973855425, 20.92033159902011, 41.65385913207093, 3.64832877264055, 8.408933464306573, 9.124965176503267, 18.14566872292388, 23.94063685625532),
    @(10.18545311809818, 5.74291052418929, 20.89026787992988, 41.63127084107355, 3.648638396664701, 8.411942808157843, 9.088104989912846, 18.11502583976509, 23.95621586053771),
    @(10.15301715422047, 5.731227283866768, 20.87190897870418, 41.61763243373337, 3.64882912802816, 8.413796059703705, 9.065468972750324, 18.09628437008288, 23.96586002196254),
    @(9.965406325157383, 5.663681163800923, 20.7675819131494, 41.5425168350291, 3.649938681253626, 8.4245689447571, 8.934895605229075, 17.98933999782004, 24.02268169414915),
    @(9.848821330815957, 5.621726914079749, 20.70439006105823, 41.49916185303933, 3.650634156410373, 8.431314270305391, 8.854068561775236, 17.92416492673398, 24.05892530103661),
    @(9.781239136888543, 5.597411390590267, 20.6683505998494, 41.47523184323443, 3.651039627588786, 8.43524428814681, 8.807328553309238, 17.88684671149664, 24.08027994633381),
    @(9.758268928086304, 5.589147400607251, 20.65620183263457, 41.46730287718903, 3.651177850854325, 8.436583572404082, 8.79146171738452, 17.87424134087477, 24.08759748147834),
    @(9.861286924058557, 5.626212239704762, 20.71108539621446, 41.50367298731707, 3.650559557903132, 8.430591017245714, 8.862699020733855, 17.93108568438083, 24.05501448327649),
    @(10.2009756342111, 5.748502261169701, 20.89908670694842, 41.63786442594608, 3.64854722806372, 8.411056819764241, 9.098943907900276, 18.12402058069667, 23.95161875730301),
    @(10.41783735869874, 5.82667013713024, 21.02448799309592, 41.73442261933106, 3.647280325397546, 8.398735369921754, 9.250785583148028, 18.25140316007824, 23.88858753409274),
    @(10.30258598185515, 5.785116348356155, 20.95733717966865, 41.68208294947743, 3.647952091891714, 8.405270937357182, 9.169994081325468, 18.18330919730467, 23.92181186285289),
    @(9.855652955601023, 5.624185031569519, 20.7080575284694, 41.50163040188592, 3.650593266337634, 8.430917837819802, 8.858798032050009, 17.92795633460749, 24.05678095149849),
    @(9.353814529434487, 5.443615333244051, 20.45100898336532, 41.34551055190084, 3.653650827952025, 8.460506250293413, 8.513774906472912, 17.65909178116968, 24.2217893843412)
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $newValues[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $rowValues[$j]
    }
}
